$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(3, "Femacal de La Calera", "Coquimbo", 45021, 5, "Fruta", 100104, "Frutos de pepita", 100104003, "Membrillo", "Champion", "Especial", 58, 16000, 16000, 16000, "$/caja 18 kilos empedrada", "Región de O'Higgins", 889, 18),
    @(3, "Femacal de La Calera", "Coquimbo", 45021, 5, "Fruta", 100104, "Frutos de pepita", 100104003, "Membrillo", "Champion", "Extra (doble especial)", 60, 18000, 18000, 18000, "$/caja 18 kilos empedrada", "Región de O'Higgins", 1000, 18),
    @(3, "Femacal de La Calera", "Coquimbo", 45021, 5, "Fruta", 100104, "Frutos de pepita", 100104003, "Membrillo", "Champion", "Primera", 57, 14000, 14000, 14000, "$/caja 18 kilos empedrada", "Región de O'Higgins", 778, 18)
)

$startRow = 109
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($c = 1; $c -le $data.Length; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($c -eq 4) {
            $cell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
        }
        $cell.Value = $data[$c - 1]
    }
}

Write-Host "Added rows 109-111"
